$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a row above the
# current row 160 (shifting the existing rows 160:251 down to 161:252)
# and populate it with the new entry's data.
$ws.Rows(160).Insert()

$ws.Range("A160").Value = 6
$ws.Range("B160").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C160").Value = "Metropolitana"
$ws.Range("D160").Value = 44873
$ws.Range("E160").Value = 13
$ws.Range("F160").Value = 100112029
$ws.Range("G160").Value = "Orégano"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 46
$ws.Range("K160").Value = 16000
$ws.Range("L160").Value = 17000
$ws.Range("M160").Value = 16457
$ws.Range("N160").Value = "`$/docena de atados"
$ws.Range("O160").Value = "Región Metropolitana"
$ws.Range("P160").Value = 5486
$ws.Range("Q160").Value = 3
$ws.Range("R160").Value = "Hortaliza"

$ws.Range("D160").NumberFormat = $ws.Range("D161").NumberFormat
